# Update cryptos list (GitHub Actions scheduled refresh).
# Prices/volumes for each coin row are refreshed; a leading "'" forces
# numeric-looking price strings to stay stored as text (matching the
# original inlineStr cells) instead of being auto-coerced to numbers by
# Excel's COM type inference. Rows 48/49 (Cronos/Aave) also swap ranking
# position, so their Coin/Link/Price/Volume cells are fully replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.638.19"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.277.23"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'112.51"
$ws.Range("D6").Value = "'266.76"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.613"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'48.39"
$ws.Range("E10").Value = "  +4.77%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'8.84"
$ws.Range("E12").Value = "  +9.71%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'15.65"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "2.619.43"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'0.872"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "2.276.62"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "43.468.10"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'7.00"
$ws.Range("E20").Value = "  +11.87%  "
$ws.Range("D21").Value = "'72.06"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D23").Value = "'9.95"
$ws.Range("E23").Value = "  +8.30%  "
$ws.Range("D24").Value = "'232.21"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'11.53"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").Value = "'41.54"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "'173.02"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "'21.49"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'0.0916"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'5.65"
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'4.62"
$ws.Range("E36").Value = "  -6.14%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").Value = "'3.73"
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("E40").Value = "  +21.38%  "
$ws.Range("D41").Value = "'74.66"
$ws.Range("E41").Value = "  +13.71%  "
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'6.27"
$ws.Range("E44").Value = "  +19.17%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'8.68"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'101.84"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0996"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "'0.453"
$ws.Range("E51").Value = "  +1.66%  "
